$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1266
$ws1.Range("F10").Value = 3502
$ws1.Range("F15").Value = 56
$ws1.Range("F18").Value = 749
$ws1.Range("F19").Value = 211
$ws1.Range("F24").Value = 2662
$ws1.Range("F25").Value = 5167
$ws1.Range("F31").Value = 2253
$ws1.Range("F35").Value = 122
$ws1.Range("F36").Value = 180
$ws1.Range("F38").Value = 24
$ws1.Range("F39").Value = 461
$ws1.Range("F40").Value = 802

# Sheet "全部类型" (sheet4) updates to column F (想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1266
$ws4.Range("F10").Value = 3502
$ws4.Range("F16").Value = 56
$ws4.Range("F19").Value = 749
$ws4.Range("F20").Value = 211
$ws4.Range("F25").Value = 2662
$ws4.Range("F26").Value = 5167
$ws4.Range("F32").Value = 2253
$ws4.Range("F36").Value = 122
$ws4.Range("F37").Value = 180
$ws4.Range("F39").Value = 24
$ws4.Range("F40").Value = 461
$ws4.Range("F41").Value = 802
